$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40: uStringArray header + array formula spilling B40:E40 ---
$ws.Range("A40").Value = "uStringArray"
$ws.Range("B40:E40").FormulaArray = "=_xll.uStringArray(B41:C42)"

# --- Rows 41-42: source data for uStringArray / uStringMatrix, right aligned ---
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = "a"
$ws.Range("B42").Value = "b"
$ws.Range("C42").Value = 2
$ws.Range("B41:C42").HorizontalAlignment = -4152

# --- Row 43-44: uStringMatrix array formula spilling B43:C44 ---
$ws.Range("B43:C44").FormulaArray = "=_xll.uStringMatrix(B41:C42)"

# --- Restore selection to the newly entered cell, as the author left it ---
[void]$ws.Range("B44").Select()
